# Refresh the cryptocurrency price/volume snapshot (GitHub Actions data pull).
#
# Column D ("Price") holds plain-text numbers (e.g. "0.998", "35.17") so that
# thousand-grouped values like "42.467.14" round-trip unchanged; assigning such
# digit-only strings straight to .Value would make Excel's COM layer coerce them
# to real numbers. Prefixing with a leading apostrophe forces text, exactly as
# typing  '0.998  into the cell by hand would.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '42.467.14'
$ws.Range('E2').Value = '  -0.74%  '

# Row 3
$ws.Range('D3').Value = '2.516.92'
$ws.Range('E3').Value = '  -1.23%  '

# Row 4
$ws.Range('D4').Value = '''0.998'
$ws.Range('E4').Value = '  -0.13%  '

# Row 5
$ws.Range('D5').Value = '''310.98'
$ws.Range('E5').Value = '  -0.24%  '

# Row 6
$ws.Range('D6').Value = '''98.86'
$ws.Range('E6').Value = '  -2.72%  '

# Row 7
$ws.Range('E7').Value = '  -1.41%  '

# Row 8
$ws.Range('E8').Value = '  -0.08%  '

# Row 9
$ws.Range('D9').Value = '''0.517'
$ws.Range('E9').Value = '  -3.18%  '

# Row 10
$ws.Range('D10').Value = '''35.17'
$ws.Range('E10').Value = '  -2.99%  '

# Row 11
$ws.Range('D11').Value = '''0.0800'
$ws.Range('E11').Value = '  -1.11%  '

# Row 12
$ws.Range('E12').Value = '  +0.16%  '

# Row 13
$ws.Range('D13').Value = '''7.20'
$ws.Range('E13').Value = '  -2.88%  '

# Row 14
$ws.Range('D14').Value = '2.903.14'
$ws.Range('E14').Value = '  -1.21%  '

# Row 15
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '2.528.06'
$ws.Range('E15').Value = '  +0.53%  '

# Row 16
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = '''15.24'
$ws.Range('E16').Value = '  -4.01%  '

# Row 17
$ws.Range('D17').Value = '''0.806'
$ws.Range('E17').Value = '  -4.04%  '

# Row 18
$ws.Range('D18').Value = '42.423.38'
$ws.Range('E18').Value = '  -0.88%  '

# Row 19
$ws.Range('D19').Value = '''6.60'
$ws.Range('E19').Value = '  -2.99%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0944'
$ws.Range('E20').Value = '  -1.36%  '

# Row 21
$ws.Range('D21').Value = '''12.03'
$ws.Range('E21').Value = '  -3.11%  '

# Row 22
$ws.Range('D22').Value = '''69.03'
$ws.Range('E22').Value = '  -0.27%  '

# Row 23
$ws.Range('D23').Value = '''240.72'
$ws.Range('E23').Value = '  -2.33%  '

# Row 24
$ws.Range('D24').Value = '''2.85'

# Row 25
$ws.Range('D25').Value = '''1.99'
$ws.Range('E25').Value = '  -3.60%  '

# Row 26
$ws.Range('E26').Value = '  +0.10%  '

# Row 27
$ws.Range('D27').Value = '''25.19'
$ws.Range('E27').Value = '  -5.27%  '

# Row 28
$ws.Range('E28').Value = '  -3.85%  '

# Row 29
$ws.Range('D29').Value = '''10.02'

# Row 30
$ws.Range('D30').Value = '''38.38'
$ws.Range('E30').Value = '  -6.65%  '

# Row 31
$ws.Range('D31').Value = '''5.86'
$ws.Range('E31').Value = '  +1.86%  '

# Row 32
$ws.Range('D32').Value = '''156.38'
$ws.Range('E32').Value = '  -0.28%  '

# Row 33
$ws.Range('D33').Value = '''2.79'
$ws.Range('E33').Value = '  +6.45%  '

# Row 34
$ws.Range('E34').Value = '  +1.49%  '

# Row 35
$ws.Range('E35').Value = '  -2.89%  '

# Row 36
$ws.Range('D36').Value = '''3.13'
$ws.Range('E36').Value = '  -4.63%  '

# Row 37
$ws.Range('E37').Value = '  -6.56%  '

# Row 38
$ws.Range('E38').Value = '  -5.08%  '

# Row 39
$ws.Range('D39').Value = '''0.107'
$ws.Range('E39').Value = '  -4.06%  '

# Row 40
$ws.Range('E40').Value = '  -1.15%  '

# Row 41
$ws.Range('D41').Value = '''4.15'
$ws.Range('E41').Value = '  -1.73%  '

# Row 42
$ws.Range('D42').Value = '''21.73'
$ws.Range('E42').Value = '  -3.19%  '

# Row 43
$ws.Range('E43').Value = '  -0.04%  '

# Row 44
$ws.Range('D44').Value = '''3.25'
$ws.Range('E44').Value = '  -1.08%  '

# Row 45
$ws.Range('E45').Value = '  -1.75%  '

# Row 46
$ws.Range('D46').Value = '1.995.99'
$ws.Range('E46').Value = '  +0.68%  '

# Row 47
$ws.Range('D47').Value = '''9.08'
$ws.Range('E47').Value = '  +0.51%  '

# Row 48
$ws.Range('D48').Value = '2.759.50'
$ws.Range('E48').Value = '  -1.20%  '

# Row 49
$ws.Range('E49').Value = '  -2.54%  '

# Row 50
$ws.Range('D50').Value = '''78.70'
$ws.Range('E50').Value = '  -3.45%  '

# Row 51
$ws.Range('D51').Value = '''100.22'
$ws.Range('E51').Value = '  -1.84%  '
